# Auto-generated Excel COM-interop edit script
# Applies the cell-level changes described in the commit diff
# (price/quantity corrections and two transposed product rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F26").Value = 317
$ws.Range("G26").Value = 16239.91
$ws.Range("B40").Value = 52198.37
$ws.Range("F44").Value = 399
$ws.Range("G44").Value = 14527.59
$ws.Range("F53").Value = 45
$ws.Range("G53").Value = 1138.05
$ws.Range("B72").Value = 170665.08
$ws.Range("B132").Value = 64196
$ws.Range("B133").Value = 65258
$ws.Range("B167").Value = 57756
$ws.Range("E167").Value = 79.37
$ws.Range("F167").Value = -100
$ws.Range("G167").Value = -6644
$ws.Range("B168").Value = 64350
$ws.Range("E168").Value = 70.63
$ws.Range("F168").Value = 2
$ws.Range("G168").Value = 132.88
$ws.Range("F190").Value = 40
$ws.Range("G190").Value = 3280
$ws.Range("B199").Value = 54546.34
$ws.Range("F223").Value = 10
$ws.Range("G223").Value = 743
$ws.Range("B224").Value = 63880.84
$ws.Range("B298").Value = 66196
$ws.Range("C298").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F298").Value = 1
$ws.Range("G298").Value = 87.7
$ws.Range("B299").Value = 64985
$ws.Range("C299").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F299").Value = 12
$ws.Range("G299").Value = 1052.4
$ws.Range("B310").Value = 63520
$ws.Range("E310").Value = 153.4
$ws.Range("F310").Value = 35
$ws.Range("G310").Value = 5049.8
$ws.Range("B311").Value = 55373
$ws.Range("E311").Value = 163.62
$ws.Range("F311").Value = -94
$ws.Range("G311").Value = -13562.32
$ws.Range("B312").Value = 63531
$ws.Range("E312").Value = 152.53
$ws.Range("F312").Value = 23
$ws.Range("G312").Value = 3300.04
$ws.Range("B313").Value = 57802
$ws.Range("E313").Value = 162.71
$ws.Range("F313").Value = -79
$ws.Range("G313").Value = -11334.92
$ws.Range("B314").Value = 63510
$ws.Range("E314").Value = 50.66
$ws.Range("F314").Value = 74
$ws.Range("G314").Value = 3525.36
$ws.Range("B315").Value = 55356
$ws.Range("E315").Value = 54.04
$ws.Range("F315").Value = -158
$ws.Range("G315").Value = -7527.12
$ws.Range("B323").Value = 60325
$ws.Range("E323").Value = 151.57
$ws.Range("F323").Value = -102
$ws.Range("G323").Value = -12939.72
$ws.Range("B324").Value = 63560
$ws.Range("E324").Value = 134.87
$ws.Range("F324").Value = 1
$ws.Range("G324").Value = 126.86
$ws.Range("F358").Value = 117
$ws.Range("G358").Value = 14614.47
$ws.Range("B362").Value = 69804.52
$ws.Range("F389").Value = 10
$ws.Range("G389").Value = 2648.3
$ws.Range("B395").Value = 49412.29
$ws.Range("F478").Value = 8
$ws.Range("G478").Value = 10142.56
$ws.Range("B482").Value = 41822.08
$ws.Range("B502").Value = 64833
$ws.Range("E502").Value = 34.9
$ws.Range("F502").Value = 88
$ws.Range("G502").Value = 2889.04
$ws.Range("B503").Value = 60025
$ws.Range("E503").Value = 37.22
$ws.Range("F503").Value = -98
$ws.Range("G503").Value = -3217.34
$ws.Range("B512").Value = 60022
$ws.Range("E512").Value = 37.22
$ws.Range("F512").Value = -113
$ws.Range("G512").Value = -3709.79
$ws.Range("B513").Value = 64830
$ws.Range("E513").Value = 34.9
$ws.Range("F513").Value = 83
$ws.Range("G513").Value = 2724.89
$ws.Range("F525").Value = 343
$ws.Range("G525").Value = 18823.84
$ws.Range("F529").Value = 133
$ws.Range("G529").Value = 11386.13
$ws.Range("B531").Value = 106729.28
$ws.Range("F680").Value = 332
$ws.Range("G680").Value = 54152.52
$ws.Range("B686").Value = 55165.07
$ws.Range("F695").Value = 22
$ws.Range("G695").Value = 2757.7
$ws.Range("B697").Value = 8897.34
$ws.Range("B724").Value = 2156014.85
$ws.Range("B725").Value = 2156014.85

Write-Host "Applied 102 cell updates"